$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the stored timestamp precision on the existing last row (A17)
$ws.Range("A17").Value = 45817.39397020833

# Append the new price entry as row 18
$ws.Range("A18").Value = 45818.39375911181
$ws.Range("A18").NumberFormat = $ws.Range("A17").NumberFormat
$ws.Range("B18").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C18").Value = "1Kg"
$ws.Range("D18").Value = "15,41€"
